# Daily attendance processing - 2026-01-03 09:34:06
# Swap the first and last comma-separated tokens in column G ("Recorded By")
# for a specific set of session rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2,3,4,5,6,7,8,28,29,30,31,32,33,34,54,55,56,57,58,59,60,80,81,82,106,107,108,132,133,134)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $value = $cell.Value2
    $parts = $value -split ", "
    if ($parts.Count -gt 1) {
        $first = $parts[0]
        $last = $parts[$parts.Count - 1]
        $parts[0] = $last
        $parts[$parts.Count - 1] = $first
        $cell.Value2 = [string]::Join(", ", $parts)
    }
}
